# Update event stats (want-to-go counts and minimum ticket prices) for
# the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

function Apply-Updates($ws, $rowOffset) {
    # $rowOffset shifts row numbers for the "全部类型" sheet, which has
    # one extra row inserted (at row 21) compared to "展览", pushing
    # everything from row 22 onward down by one row.

    $ws.Range("F3").Value = 5525
    $ws.Range("F4").Value = 46
    $ws.Range("G6").Value = 60
    $ws.Range("G8").Value = 50
    $ws.Range("F9").Value = 4371
    $ws.Range("G9").Value = 80
    $ws.Range("G10").Value = 55
    $ws.Range("F11").Value = 814
    $ws.Range("G11").Value = 65
    $ws.Range("G12").Value = 55
    $ws.Range("F14").Value = 126
    $ws.Range("F15").Value = 143
    $ws.Range("G15").Value = "不可售"
    $ws.Range("F18").Value = 123
    $ws.Range("F20").Value = 20

    $r22 = 22 + $rowOffset
    $r23 = 23 + $rowOffset
    $r24 = 24 + $rowOffset
    $r25 = 25 + $rowOffset
    $r26 = 26 + $rowOffset

    $ws.Range("F$r22").Value = 1133
    $ws.Range("F$r23").Value = 16
    $ws.Range("F$r24").Value = 2762
    $ws.Range("F$r25").Value = 440
    $ws.Range("F$r26").Value = 292
}

$wsExhibit = $wb.Worksheets.Item("展览")
Apply-Updates $wsExhibit 0

$wsAll = $wb.Worksheets.Item("全部类型")
Apply-Updates $wsAll 1
